$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.703.44'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.26%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.476.82'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.74%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '320.20'
$ws.Range('D5').Style = 'Normal'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '92.32'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.41%  '

$ws.Range('E7').Value = '  +0.80%  '

$ws.Range('E8').Value = '  +0.04%  '

$ws.Range('E9').Value = '  +0.37%  '

$ws.Range('E10').Value = '  +5.71%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '33.04'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.13%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.109'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.78%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.858.18'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.68%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.90'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.96%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.49'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.84%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.460.67'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.52%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.794'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.34%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '41.621.79'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.03%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.44'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.34%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0944'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.21%  '

$ws.Range('E21').Value = '  +0.11%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.26'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.28%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '239.67'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.81%  '

$ws.Range('E24').Value = '  +1.66%  '

$ws.Range('E25').Value = '  +2.42%  '

$ws.Range('E26').Value = '  +0.03%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.23%  '

$ws.Range('E28').Value = '  -0.60%  '

$ws.Range('E29').Value = '  +0.67%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.69'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.86%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '157.60'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.25%  '

$ws.Range('E32').Value = '  -0.38%  '

$ws.Range('E33').Value = '  -0.03%  '

$ws.Range('E34').Value = '  +1.04%  '

$ws.Range('E35').Value = '  -0.03%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.18'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.28%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.84'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.03%  '

$ws.Range('E38').Value = '  +1.66%  '

$ws.Range('E39').Value = '  +0.29%  '

$ws.Range('E40').Value = '  +1.18%  '

$ws.Range('E41').Value = '  +2.74%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.43'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.98%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.996.84'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.04%  '

$ws.Range('E44').Value = '  +0.72%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.68'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.83%  '

$ws.Range('E46').Value = '  +2.28%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.49'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.58%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.753.67'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.09%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '97.56'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.85%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '75.99'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.70%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '67.49'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.20%  '
